$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.250615119934082
$ws.Range("B1").Value = 2.710794448852539
$ws.Range("C1").Value = 2.446203708648682
$ws.Range("D1").Value = 2.731816530227661
$ws.Range("E1").Value = 3.199106693267822
